$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card9")

# Update existing row 27: columns B..K and M become the literal string "nan"
# (they used to be blank, the other rows on this sheet already use "nan" there)
$ws.Range("B27:K27").Value = "nan"
$ws.Range("M27").Value = "nan"

# Add new row 28 - new service event for Card9 (B28:K28 stay blank, same as row 27 was)
# A28 and M28 hold numeric-looking text, so force the Text format first so Excel
# keeps them as strings (matching every other "card"/"Event" cell in this sheet)
# instead of silently converting them to numbers.
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = "9"
$ws.Range("L28").Value = "26\11\2025"
$ws.Range("M28").NumberFormat = "@"
$ws.Range("M28").Value = "1013.6"
$ws.Range("N28").Value = "تم سن الفلاتس"
$ws.Range("O28").Value = "الخبير"
